$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for "بلاستر مترسيلك 2 سم" (row 87). This shifts all
# following rows up by one and Excel automatically fixes merged cell ranges.
$ws.Rows.Item(87).Delete()

# After the deletion, the "سرنجات 3 سم" row (previously row 91) is now row 90.
# Update its sale price and transaction-count columns.
$ws.Cells.Item(90, 16).Value = "50.0000"
$ws.Cells.Item(90, 17).Value = "25:0"

# Update the generated-on timestamp in the footer (previously row 108, now row 107).
$ws.Cells.Item(107, 1).Value = "Thursday, 9 October, 2025 9:04 PM"
